# Update shooter LUT and pickup G path
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New Velocity (column B) / Hood angle (column C) values for rows 2-9
$values = @(
    @(5.603, 9250),
    @(7.298, 9750),
    @(7.9344, 10650),
    @(10.342, 11800),
    @(11.598, 12600),
    @(12.556, 11800),
    @(14.174, 12250),
    @(16.925, 13100)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i][0]
    $ws.Cells.Item($row, 3).Value = $values[$i][1]
}

# Update the selection to match the recorded view state (A3:C3, active cell A3)
$ws.Range("A3:C3").Select()
